$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.341.61'
$ws.Range('E2').Value = '  +0.85%  '

$ws.Range('D3').Value = '3.919.06'
$ws.Range('E3').Value = '  -1.46%  '

$ws.Range('E4').Value = '  +0.30%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.99%  '

$ws.Range('E7').Value = '  -0.65%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('E9').Value = '  -0.08%  '

$ws.Range('E10').Value = '  -2.27%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000346'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.45%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.56%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.77'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.74%  '

$ws.Range('D14').Value = '4.541.75'
$ws.Range('E14').Value = '  -1.42%  '

$ws.Range('D15').Value = '3.929.35'
$ws.Range('E15').Value = '  -1.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.76%  '

$ws.Range('E17').Value = '  -1.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.31%  '

$ws.Range('E19').Value = '  -1.25%  '

$ws.Range('D20').Value = '68.425.31'
$ws.Range('E20').Value = '  +0.91%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.92%  '

$ws.Range('E22').Value = '  +4.77%  '

$ws.Range('E23').Value = '  +1.65%  '

$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +17.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.85%  '

$ws.Range('E27').Value = '  +0.22%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.47%  '

$ws.Range('E29').Value = '  -0.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '714.83'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.51%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.94%  '

$ws.Range('E32').Value = '  -0.52%  '

$ws.Range('E33').Value = '  +4.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.18'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.74%  '

$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.75'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.12%  '

$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0884'
$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '60.85'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.78%  '

$ws.Range('E38').Value = '  -0.06%  '

$ws.Range('E39').Value = '  -5.47%  '

$ws.Range('E40').Value = '  +15.55%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0498'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.30%  '

$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.90%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.11'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.78%  '

$ws.Range('E44').Value = '  +4.68%  '

$ws.Range('E45').Value = '  +5.53%  '

$ws.Range('E46').Value = '  -0.84%  '

$ws.Range('E47').Value = '  +0.25%  '

$ws.Range('E48').Value = '  -1.31%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.17%  '

$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.41%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0342'
$ws.Range('E51').Value = '  +26.39%  '
